$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.018227577209473
$ws.Range("B1").Value = 2.30049729347229
$ws.Range("C1").Value = 7.886225700378418
$ws.Range("D1").Value = 1.009856224060059
$ws.Range("E1").Value = 0.5756622552871704
